# Applies the "Deals" sheet population edit described by the diff:
#  - Renames Sheet3 -> Deals
#  - Fills in the Deals header row + two data rows
#  - Applies the same look (yellow header fill, quote-prefixed numeric text)
#    used on the other sheets in this workbook
#  - Updates the selection on Companies (A1:D1) and Deals (F11)

$wb = $excel.ActiveWorkbook

# --- Deals sheet (formerly "Sheet3") ---------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Name = "Deals"

# Header row (typed left to right)
$ws3.Range("A1").Value = "title"
$ws3.Range("B1").Value = "clientLookup"
$ws3.Range("C1").Value = "amount"
$ws3.Range("D1").Value = "probability"
$ws3.Range("E1").Value = "commission"
$ws3.Range("F1").Value = "identifier"
$ws3.Range("G1").Value = "quantity"
$ws3.Range("H1").Value = "type"
$ws3.Range("I1").Value = "source"
$ws3.Range("A1:I1").Interior.Color = 65535

# Data rows, filled column by column (matches original authoring order)
$ws3.Range("A2").Value = "New Deal 1"
$ws3.Range("A3").Value = "New Deal 2"

$ws3.Range("B2").Value = "Amazon"
$ws3.Range("B3").Value = "Oracle"

$ws3.Range("C2").Value = "'10000000"
$ws3.Range("C3").Value = "'200000000"

$ws3.Range("D2").Value = "'100"
$ws3.Range("D3").Value = "'80"

$ws3.Range("E2").Value = "'10"
$ws3.Range("E3").Value = "'12"

$ws3.Range("F2").Value = "User1"
$ws3.Range("F3").Value = "User2"

$ws3.Range("G2").Value = "'10000"

$ws3.Range("H2").Value = "New"
$ws3.Range("H3").Value = "Opportunity"

$ws3.Range("I2").Value = "Existing Customer"
$ws3.Range("I3").Value = "Partner"

$ws3.Range("G3").Value = "'5000"

$ws3.Columns.Item(1).EntireColumn.AutoFit()
$ws3.Columns.Item(2).EntireColumn.AutoFit()
$ws3.Columns.Item(3).EntireColumn.AutoFit()
$ws3.Columns.Item(4).EntireColumn.AutoFit()
$ws3.Columns.Item(5).EntireColumn.AutoFit()
$ws3.Columns.Item(6).EntireColumn.AutoFit()
$ws3.Columns.Item(7).EntireColumn.AutoFit()
$ws3.Columns.Item(8).EntireColumn.AutoFit()
$ws3.Columns.Item(9).EntireColumn.AutoFit()

$null = $ws3.Range("F11").Select()

# --- Companies sheet: remains the active tab, selection changes to A1:D1 ---
$ws2 = $wb.Worksheets.Item("Companies")
$null = $ws2.Activate()
$null = $ws2.Range("A1:D1").Select()
